$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.135083983912941
$ws.Range("D2").Value = "3rd Tier"
$ws.Range("C3").Value = 0.6841731724627396
$ws.Range("D3").Value = "Below Median"
$ws.Range("C4").Value = 1.685592618878637
$ws.Range("D4").Value = "1st Tier"
$ws.Range("C5").Value = 1.376389874615566
$ws.Range("D5").Value = "2nd Tier"
$ws.Range("C6").Value = 1.685592618878637
$ws.Range("D6").Value = "1st Tier"
$ws.Range("C7").Value = 2.045422285308729
$ws.Range("D7").Value = "1st Tier"
$ws.Range("C8").Value = 1.325762952448545
$ws.Range("D8").Value = "2nd Tier"
$ws.Range("C9").Value = 2.045422285308729
$ws.Range("D9").Value = "1st Tier"
$ws.Range("C10").Value = 0.3553347527797492
$ws.Range("D10").Value = "Below Median"
$ws.Range("C11").Value = 0.4178850248403123
$ws.Range("D11").Value = "Below Median"
$ws.Range("C12").Value = 0.6728176011355571
$ws.Range("D12").Value = "Below Median"
$ws.Range("C13").Value = 0.9609652235628106
$ws.Range("D13").Value = "Below Median"
$ws.Range("C14").Value = 1.163946061036196
$ws.Range("D14").Value = "3rd Tier"
$ws.Range("C15").Value = 1.163946061036196
$ws.Range("D15").Value = "3rd Tier"
$ws.Range("C16").Value = 0.7324343506032647
$ws.Range("D16").Value = "Below Median"
$ws.Range("C17").Value = 0.9247693399574166
$ws.Range("D17").Value = "Below Median"
$ws.Range("C18").Value = 0.7686302342086586
$ws.Range("D18").Value = "Below Median"
$ws.Range("C19").Value = 0.8833688194937308
$ws.Range("D19").Value = "Below Median"
$ws.Range("C20").Value = 0.8833688194937308
$ws.Range("D20").Value = "Below Median"
$ws.Range("C21").Value = 1.185237757274663
$ws.Range("D21").Value = "3rd Tier"
$ws.Range("C22").Value = 0.5613910574875799
$ws.Range("D22").Value = "Below Median"
$ws.Range("C23").Value = 0.5613910574875799
$ws.Range("D23").Value = "Below Median"
$ws.Range("C24").Value = 1.182398864442867
$ws.Range("D24").Value = "3rd Tier"
$ws.Range("C25").Value = 1.611071682044003
$ws.Range("D25").Value = "2nd Tier"
$ws.Range("C26").Value = 1.753726046841732
$ws.Range("D26").Value = "1st Tier"
$ws.Range("C27").Value = 1.753726046841732
$ws.Range("D27").Value = "1st Tier"
$ws.Range("C28").Value = 1.474095102909865
$ws.Range("D28").Value = "2nd Tier"
$ws.Range("C29").Value = 0.8211497515968772
$ws.Range("D29").Value = "Below Median"
$ws.Range("C30").Value = 1.249112845990064
$ws.Range("D30").Value = "2nd Tier"
$ws.Range("C31").Value = 0.9141234918381831
$ws.Range("D31").Value = "Below Median"
$ws.Range("C32").Value = 1.182398864442867
$ws.Range("D32").Value = "3rd Tier"
$ws.Range("C33").Value = 1.039744499645138
$ws.Range("D33").Value = "4th Tier"
$ws.Range("C34").Value = 0.9921930447125621
$ws.Range("D34").Value = "Below Median"
$ws.Range("C35").Value = 0.7314880529926663
$ws.Range("D35").Value = "Below Median"
$ws.Range("C36").Value = 1.028152353915306
$ws.Range("D36").Value = "4th Tier"
$ws.Range("C37").Value = 0.9027679205110007
$ws.Range("D37").Value = "Below Median"
$ws.Range("C38").Value = 0.9027679205110007
$ws.Range("D38").Value = "Below Median"
$ws.Range("C39").Value = 1.462029808374734
$ws.Range("D39").Value = "2nd Tier"
$ws.Range("C40").Value = 1.028152353915306
$ws.Range("D40").Value = "4th Tier"
$ws.Range("C41").Value = 0.9027679205110007
$ws.Range("D41").Value = "Below Median"
$ws.Range("C42").Value = 0.9623846699787083
$ws.Range("D42").Value = "Below Median"
$ws.Range("C43").Value = 1.182398864442867
$ws.Range("D43").Value = "3rd Tier"
$ws.Range("C44").Value = 1.12230896616986
$ws.Range("D44").Value = "4th Tier"
$ws.Range("C45").Value = 1.135083983912941
$ws.Range("D45").Value = "3rd Tier"
$ws.Range("C46").Value = 0.6898509581263307
$ws.Range("D46").Value = "Below Median"
$ws.Range("C47").Value = 1.039744499645138
$ws.Range("D47").Value = "4th Tier"
$ws.Range("C48").Value = 0.890347764371895
$ws.Range("D48").Value = "Below Median"
$ws.Range("C49").Value = 0.7885024840312278
$ws.Range("D49").Value = "Below Median"
$ws.Range("C50").Value = 0.7314880529926663
$ws.Range("D50").Value = "Below Median"
$ws.Range("C51").Value = 1.135083983912941
$ws.Range("D51").Value = "3rd Tier"
$ws.Range("C52").Value = 1.039744499645138
$ws.Range("D52").Value = "4th Tier"
$ws.Range("C53").Value = 1.090844570617459
$ws.Range("D53").Value = "4th Tier"
$ws.Range("C54").Value = 0.4059616749467708
$ws.Range("D54").Value = "Below Median"
$ws.Range("C55").Value = 0.9854506742370476
$ws.Range("D55").Value = "Below Median"
$ws.Range("C56").Value = 1.235391530636385
$ws.Range("D56").Value = "3rd Tier"
$ws.Range("C57").Value = 0.5734563520227112
$ws.Range("D57").Value = "Below Median"
$ws.Range("C58").Value = 0.5956943458717767
$ws.Range("D58").Value = "Below Median"
$ws.Range("C59").Value = 0.9666430092264017
$ws.Range("D59").Value = "Below Median"
$ws.Range("C60").Value = 0.850958126330731
$ws.Range("D60").Value = "Below Median"
$ws.Range("C61").Value = 1.325762952448545
$ws.Range("D61").Value = "2nd Tier"
$ws.Range("C62").Value = 1.462029808374734
$ws.Range("D62").Value = "2nd Tier"
$ws.Range("C63").Value = 0.9623846699787083
$ws.Range("D63").Value = "Below Median"
$ws.Range("C64").Value = 0.5956943458717767
$ws.Range("D64").Value = "Below Median"
$ws.Range("C65").Value = 1.039744499645138
$ws.Range("D65").Value = "4th Tier"
$ws.Range("C66").Value = 0.8034066713981547
$ws.Range("D66").Value = "Below Median"
$ws.Range("C67").Value = 1.082564466524722
$ws.Range("D67").Value = "4th Tier"
$ws.Range("C68").Value = 1.611071682044003
$ws.Range("D68").Value = "1st Tier"
$ws.Range("C69").Value = 0.3222143364088006
$ws.Range("D69").Value = "Below Median"
$ws.Range("C70").Value = 0.6174591909155429
$ws.Range("D70").Value = "Below Median"
$ws.Range("C71").Value = 1.611071682044003
$ws.Range("D71").Value = "1st Tier"
$ws.Range("C72").Value = 0.9609652235628106
$ws.Range("D72").Value = "Below Median"
$ws.Range("C73").Value = 0.4932576295244854
$ws.Range("D73").Value = "Below Median"
$ws.Range("C74").Value = 1.323278921220724
$ws.Range("D74").Value = "2nd Tier"
$ws.Range("C75").Value = 1.685592618878637
$ws.Range("D75").Value = "1st Tier"
$ws.Range("C76").Value = 1.611071682044003
$ws.Range("D76").Value = "1st Tier"
$ws.Range("C77").Value = 1
$ws.Range("D77").Value = "4th Tier"
$ws.Range("C78").Value = 1.470073338064821
$ws.Range("D78").Value = "2nd Tier"
$ws.Range("C79").Value = 0.971611071682044
$ws.Range("D79").Value = "Below Median"
$ws.Range("C80").Value = 0.4758694109297374
$ws.Range("D80").Value = "Below Median"
